# Updated symbol list on Sun Dec 18 03:30:24 UTC 2022 with GitHub Actions
# Refreshes the Price (column D) and Worst/Best-in-24h label (column E)
# cells for the symbols whose market data changed in this run.
#
# Column D values are numeric-looking text (e.g. "241.88", "0.6800") that
# must stay stored as literal text (so trailing zeros / exact formatting
# survive) rather than being coerced to a Number -- hence the leading
# apostrophe, which is the standard Excel way to force text entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'241.88"
$ws.Range("D4").Value  = "'5.522"
$ws.Range("D5").Value  = "'0.05585"
$ws.Range("D6").Value  = "'3.383"
$ws.Range("D7").Value  = "'6.471"
$ws.Range("D8").Value  = "'1.079"
$ws.Range("D9").Value  = "'0.8030"
$ws.Range("D10").Value = "'0.1418"
$ws.Range("D11").Value = "'0.07428"
$ws.Range("D12").Value = "'0.03258"
$ws.Range("D13").Value = "'0.02980"
$ws.Range("D14").Value = "'0.09254"
$ws.Range("D15").Value = "'0.001675"
$ws.Range("D16").Value = "'3.253"
$ws.Range("D17").Value = "'0.04713"
$ws.Range("D18").Value = "'0.0005740"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006262"
$ws.Range("D20").Value = "'0.001047"
$ws.Range("D21").Value = "'0.003802"
$ws.Range("E21").Value = "20HotbitTokenHTBWorstin24h"
$ws.Range("D24").Value = "'3.981"
$ws.Range("D27").Value = "'0.1312"
$ws.Range("D40").Value = "'0.04175"
$ws.Range("D41").Value = "'0.007060"
$ws.Range("D42").Value = "'0.003500"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("D44").Value = "'0.009021"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("D45").Value = "'0.00005496"
$ws.Range("D47").Value = "'0.6800"
$ws.Range("D48").Value = "'0.03032"
